$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A68").Value = "Human Resources  Position Catagories Add Functionality"
$ws.Range("B68").Value = "PASSED"
$ws.Range("C68").Value = "chrome"

$ws.Range("A69").Value = "Human Resources  PositionCatagories Edit Functionality"
$ws.Range("B69").Value = "PASSED"
$ws.Range("C69").Value = "chrome"

$ws.Range("A70").Value = "Human Resources  PositionCatagories Delete Functionality"
$ws.Range("B70").Value = "PASSED"
$ws.Range("C70").Value = "chrome"
